$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 335572.34
$ws.Range("J17").Value = 359398.97
$ws.Range("L17").Value = 1078196.91
$ws.Range("N17").Value = -1078532.91
$ws.Range("H33").Value = 701.9167
$ws.Range("I33").Value = 255
$ws.Range("K33").Value = 255
$ws.Range("M33").Value = -26
$ws.Range("H53").Value = 1740.75
$ws.Range("I53").Value = 1912.8572
$ws.Range("J53").Value = 1499.8
$ws.Range("K53").Value = 1912.8572
$ws.Range("L53").Value = 1499.8
$ws.Range("M53").Value = -1275.8572
$ws.Range("N53").Value = -2773.8
$ws.Range("H87").Value = 67584.75
$ws.Range("J87").Value = 67584.75
$ws.Range("L87").Value = 67584.75
$ws.Range("N87").Value = -70080.75
$ws.Range("H90").Value = 67584.75
$ws.Range("J90").Value = 67584.75
$ws.Range("L90").Value = 202754.25
$ws.Range("N90").Value = -215234.25
$ws.Range("H92").Value = 3430.0557
$ws.Range("I92").Value = 3724.913
$ws.Range("J92").Value = 2908.3845
$ws.Range("K92").Value = 3724.913
$ws.Range("L92").Value = 2908.3845
$ws.Range("M92").Value = -2476.913
$ws.Range("N92").Value = -5404.3845
$ws.Range("H95").Value = 31500
$ws.Range("J95").Value = 31500
$ws.Range("L95").Value = 31500
$ws.Range("N95").Value = -36992
$ws.Range("H106").Value = 130402.875
$ws.Range("I106").Value = 130402.875
$ws.Range("K106").Value = 130402.875
$ws.Range("M106").Value = -129771.875
$ws.Range("H116").Value = 10929.143
$ws.Range("J116").Value = 5848.8
$ws.Range("L116").Value = 5848.8
$ws.Range("N116").Value = -12732.8
$ws.Range("H118").Value = 125000230
$ws.Range("I118").Value = 166666820
$ws.Range("J118").Value = 499
$ws.Range("K118").Value = 500000460
$ws.Range("L118").Value = 1497
$ws.Range("M118").Value = -499998803
$ws.Range("N118").Value = -4811
$ws.Range("H131").Value = 21134.092
$ws.Range("I131").Value = 21134.092
$ws.Range("K131").Value = 63402.276
$ws.Range("M131").Value = -58362.276
$ws.Range("H132").Value = 1619.2273
$ws.Range("I132").Value = 1438.3636
$ws.Range("J132").Value = 2161.818
$ws.Range("K132").Value = 4315.0908
$ws.Range("L132").Value = 6485.454000000001
$ws.Range("M132").Value = -1785.0908
$ws.Range("N132").Value = -11545.454
$ws.Range("H137").Value = 4598141.5
$ws.Range("I137").Value = 10112491
$ws.Range("J137").Value = 2850.1667
$ws.Range("K137").Value = 30337473
$ws.Range("L137").Value = 8550.500100000001
$ws.Range("M137").Value = -30334923
$ws.Range("N137").Value = -13650.5001
$ws.Range("H138").Value = 3991.79
$ws.Range("I138").Value = 1775.3684
$ws.Range("J138").Value = 4511.6914
$ws.Range("K138").Value = 5326.1052
$ws.Range("L138").Value = 13535.0742
$ws.Range("M138").Value = -186.1052
$ws.Range("N138").Value = -23815.0742

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5141.625
$ws.Range("I32").Value = 5141.625
$ws.Range("K32").Value = 5141.625
$ws.Range("M32").Value = -4854.625
$ws.Range("H74").Value = 154011
$ws.Range("I74").Value = 181749
$ws.Range("J74").Value = 1452
$ws.Range("K74").Value = 181749
$ws.Range("L74").Value = 1452
$ws.Range("M74").Value = -180875
$ws.Range("N74").Value = -3200
$ws.Range("H77").Value = 154011
$ws.Range("I77").Value = 181749
$ws.Range("J77").Value = 1452
$ws.Range("K77").Value = 908745
$ws.Range("L77").Value = 7260
$ws.Range("M77").Value = -904377
$ws.Range("N77").Value = -15996
$ws.Range("H97").Value = 981.43475
$ws.Range("I97").Value = 574.05554
$ws.Range("J97").Value = 2448
$ws.Range("K97").Value = 574.05554
$ws.Range("L97").Value = 2448
$ws.Range("M97").Value = -78.05553999999995
$ws.Range("N97").Value = -3440
$ws.Range("H101").Value = 94035.57000000001
$ws.Range("J101").Value = 94035.57000000001
$ws.Range("L101").Value = 94035.57000000001
$ws.Range("N101").Value = -100525.57
$ws.Range("H110").Value = 4433.857
$ws.Range("I110").Value = 4047.8333
$ws.Range("K110").Value = 4047.8333
$ws.Range("M110").Value = -2002.8333
$ws.Range("H122").Value = 1914.5
$ws.Range("I122").Value = 1914.5
$ws.Range("K122").Value = 5743.5
$ws.Range("M122").Value = -3293.5
$ws.Range("H132").Value = 2930
$ws.Range("I132").Value = 2767.2
$ws.Range("J132").Value = 7000
$ws.Range("K132").Value = 8301.599999999999
$ws.Range("L132").Value = 21000
$ws.Range("M132").Value = -5771.599999999999
$ws.Range("N132").Value = -26060

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4526.143
$ws.Range("I20").Value = 4134.8
$ws.Range("J20").Value = 5504.5
$ws.Range("K20").Value = 4134.8
$ws.Range("L20").Value = 5504.5
$ws.Range("M20").Value = -3887.8
$ws.Range("N20").Value = -5998.5
$ws.Range("H86").Value = 3019
$ws.Range("I86").Value = 2083.8333
$ws.Range("K86").Value = 2083.8333
$ws.Range("M86").Value = -960.8332999999998
$ws.Range("H89").Value = 3019
$ws.Range("I89").Value = 2083.8333
$ws.Range("K89").Value = 10419.1665
$ws.Range("M89").Value = -4803.166499999999
$ws.Range("H99").Value = 2636.8298
$ws.Range("I99").Value = 2164.0286
$ws.Range("K99").Value = 2164.0286
$ws.Range("M99").Value = -666.0286000000001
$ws.Range("H105").Value = 1836.1578
$ws.Range("I105").Value = 1863.7858
$ws.Range("J105").Value = 1758.8
$ws.Range("K105").Value = 1863.7858
$ws.Range("L105").Value = 1758.8
$ws.Range("M105").Value = -116.7858000000001
$ws.Range("N105").Value = -5252.8
$ws.Range("H134").Value = 3790.25
$ws.Range("I134").Value = 3388.6667
$ws.Range("K134").Value = 10166.0001
$ws.Range("M134").Value = -7631.000100000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2375.1155
$ws.Range("I16").Value = 2400.12
$ws.Range("K16").Value = 2400.12
$ws.Range("M16").Value = -2113.12
$ws.Range("H31").Value = 175826.42
$ws.Range("I31").Value = 240035.78
$ws.Range("K31").Value = 240035.78
$ws.Range("M31").Value = -239740.78
$ws.Range("H34").Value = 175826.42
$ws.Range("I34").Value = 240035.78
$ws.Range("K34").Value = 240035.78
$ws.Range("M34").Value = -239833.78
$ws.Range("H50").Value = 15000
$ws.Range("J50").Value = 15000
$ws.Range("L50").Value = 15000
$ws.Range("N50").Value = -16250
$ws.Range("H55").Value = 6326.3335
$ws.Range("I55").Value = 4999
$ws.Range("K55").Value = 4999
$ws.Range("M55").Value = -4684
$ws.Range("H58").Value = 2918.25
$ws.Range("I58").Value = 2602.2222
$ws.Range("J58").Value = 3866.3333
$ws.Range("K58").Value = 2602.2222
$ws.Range("L58").Value = 3866.3333
$ws.Range("M58").Value = -2399.2222
$ws.Range("N58").Value = -4272.3333
$ws.Range("H105").Value = 4336.125
$ws.Range("I105").Value = 1674.2632
$ws.Range("K105").Value = 1674.2632
$ws.Range("M105").Value = 72.7367999999999
$ws.Range("H107").Value = 3425.1875
$ws.Range("J107").Value = 5171.933
$ws.Range("L107").Value = 5171.933
$ws.Range("N107").Value = -9011.933000000001
$ws.Range("H113").Value = 2375.1155
$ws.Range("I113").Value = 2400.12
$ws.Range("K113").Value = 2400.12
$ws.Range("M113").Value = -230.1199999999999
$ws.Range("H122").Value = 1000
$ws.Range("I122").Value = 1000
$ws.Range("K122").Value = 3000
$ws.Range("M122").Value = -550
$ws.Range("H132").Value = 3357.7693
$ws.Range("I132").Value = 3053.0833
$ws.Range("J132").Value = 7014
$ws.Range("K132").Value = 9159.249899999999
$ws.Range("L132").Value = 21042
$ws.Range("M132").Value = -6629.249899999999
$ws.Range("N132").Value = -26102
$ws.Range("H135").Value = 164999.33
$ws.Range("J135").Value = 164999.33
$ws.Range("L135").Value = 164999.33
$ws.Range("N135").Value = -175139.33
$ws.Range("H136").Value = 2918.25
$ws.Range("I136").Value = 2602.2222
$ws.Range("J136").Value = 3866.3333
$ws.Range("K136").Value = 7806.6666
$ws.Range("L136").Value = 11598.9999
$ws.Range("M136").Value = -5256.6666
$ws.Range("N136").Value = -16698.9999
$ws.Range("H141").Value = 191527.33
$ws.Range("J141").Value = 191527.33
$ws.Range("L141").Value = 191527.33
$ws.Range("N141").Value = -201887.33

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 288.41177
$ws.Range("I12").Value = 288
$ws.Range("J12").Value = 288.58334
$ws.Range("K12").Value = 864
$ws.Range("L12").Value = 865.7500200000001
$ws.Range("M12").Value = -691
$ws.Range("N12").Value = -1211.75002
$ws.Range("H36").Value = 1724
$ws.Range("I36").Value = 1724
$ws.Range("K36").Value = 5172
$ws.Range("M36").Value = -5003
$ws.Range("H39").Value = 1493.5278
$ws.Range("J39").Value = 1829.5217
$ws.Range("L39").Value = 5488.5651
$ws.Range("N39").Value = -6076.5651
$ws.Range("H46").Value = 6192193
$ws.Range("I46").Value = 4556
$ws.Range("K46").Value = 13668
$ws.Range("M46").Value = -13577
$ws.Range("H47").Value = 12079.167
$ws.Range("I47").Value = 7500.7144
$ws.Range("K47").Value = 22502.1432
$ws.Range("M47").Value = -22071.1432
$ws.Range("H55").Value = 7693.846
$ws.Range("J55").Value = 8002.4
$ws.Range("L55").Value = 24007.2
$ws.Range("N55").Value = -24361.2
$ws.Range("H60").Value = 1964.5
$ws.Range("I60").Value = 1602.2858
$ws.Range("J60").Value = 4500
$ws.Range("K60").Value = 4806.857400000001
$ws.Range("L60").Value = 13500
$ws.Range("M60").Value = -4555.857400000001
$ws.Range("N60").Value = -14002
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H114").Value = 721.7273
$ws.Range("I114").Value = 576.5
$ws.Range("K114").Value = 1729.5
$ws.Range("M114").Value = 1524.5
$ws.Range("H129").Value = 39143116
$ws.Range("H131").Value = 1440.01
$ws.Range("I131").Value = 1299
$ws.Range("J131").Value = 1441.4343
$ws.Range("K131").Value = 3897
$ws.Range("L131").Value = 4324.3029
$ws.Range("M131").Value = 1143
$ws.Range("N131").Value = -14404.3029
$ws.Range("H132").Value = 4269
$ws.Range("I132").Value = 5015.3228
$ws.Range("J132").Value = 1377
$ws.Range("K132").Value = 45137.9052
$ws.Range("L132").Value = 12393
$ws.Range("M132").Value = -42607.9052
$ws.Range("N132").Value = -17453
$ws.Range("H137").Value = 19610788
$ws.Range("I137").Value = 1397
$ws.Range("J137").Value = 41671356
$ws.Range("K137").Value = 4191
$ws.Range("L137").Value = 125014068
$ws.Range("M137").Value = 909
$ws.Range("N137").Value = -125024268
$ws.Range("H141").Value = 2575.6
$ws.Range("I141").Value = 2575.6
$ws.Range("K141").Value = 7726.799999999999
$ws.Range("M141").Value = -2546.799999999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 4
$ws.Range("I5").Value = 4
$ws.Range("K5").Value = 4
$ws.Range("M5").Value = 108
$ws.Range("H64").Value = 80000
$ws.Range("J64").Value = 80000
$ws.Range("L64").Value = 80000
$ws.Range("N64").Value = -80496
$ws.Range("H67").Value = 80000
$ws.Range("J67").Value = 80000
$ws.Range("L67").Value = 80000
$ws.Range("N67").Value = -81716
$ws.Range("H80").Value = 4711.579
$ws.Range("I80").Value = 4626.6875
$ws.Range("J80").Value = 5164.3335
$ws.Range("K80").Value = 4626.6875
$ws.Range("L80").Value = 5164.3335
$ws.Range("M80").Value = -3628.6875
$ws.Range("N80").Value = -7160.3335
$ws.Range("H83").Value = 4711.579
$ws.Range("I83").Value = 4626.6875
$ws.Range("J83").Value = 5164.3335
$ws.Range("K83").Value = 23133.4375
$ws.Range("L83").Value = 25821.6675
$ws.Range("M83").Value = -18141.4375
$ws.Range("N83").Value = -35805.6675
$ws.Range("H97").Value = 1878.4615
$ws.Range("I97").Value = 1530.3125
$ws.Range("J97").Value = 2435.5
$ws.Range("K97").Value = 1530.3125
$ws.Range("L97").Value = 2435.5
$ws.Range("M97").Value = -1034.3125
$ws.Range("N97").Value = -3427.5
$ws.Range("H102").Value = 1498.3182
$ws.Range("I102").Value = 1403.4
$ws.Range("K102").Value = 1403.4
$ws.Range("M102").Value = 218.5999999999999
$ws.Range("H113").Value = 3319.9092
$ws.Range("I113").Value = 2862.0715
$ws.Range("J113").Value = 4121.125
$ws.Range("K113").Value = 2862.0715
$ws.Range("L113").Value = 4121.125
$ws.Range("M113").Value = -692.0715
$ws.Range("N113").Value = -8461.125
$ws.Range("H122").Value = 2738.4
$ws.Range("I122").Value = 3149.5
$ws.Range("J122").Value = 1094
$ws.Range("K122").Value = 9448.5
$ws.Range("L122").Value = 3282
$ws.Range("M122").Value = -6998.5
$ws.Range("N122").Value = -8182
$ws.Range("H126").Value = 6143.5713
$ws.Range("J126").Value = 4139
$ws.Range("L126").Value = 12417
$ws.Range("N126").Value = -17357
$ws.Range("H132").Value = 59343.168
$ws.Range("I132").Value = 61466.78
$ws.Range("K132").Value = 184400.34
$ws.Range("M132").Value = -181870.34
$ws.Range("H136").Value = 8644.0625
$ws.Range("J136").Value = 8644.0625
$ws.Range("L136").Value = 25932.1875
$ws.Range("N136").Value = -31032.1875

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8749.154
$ws.Range("I7").Value = 9118.849
$ws.Range("J7").Value = 6715.8335
$ws.Range("K7").Value = 9118.849
$ws.Range("L7").Value = 6715.8335
$ws.Range("M7").Value = -9006.849
$ws.Range("N7").Value = -6939.8335
$ws.Range("H22").Value = 936.8333
$ws.Range("I22").Value = 955.25
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 955.25
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = -660.25
$ws.Range("N22").Value = -1490
$ws.Range("H27").Value = 936.8333
$ws.Range("I27").Value = 955.25
$ws.Range("J27").Value = 900
$ws.Range("K27").Value = 955.25
$ws.Range("L27").Value = 900
$ws.Range("M27").Value = -848.25
$ws.Range("N27").Value = -1114
$ws.Range("H40").Value = 15861.286
$ws.Range("I40").Value = 16159.35
$ws.Range("J40").Value = 9900
$ws.Range("K40").Value = 16159.35
$ws.Range("L40").Value = 9900
$ws.Range("M40").Value = -16023.35
$ws.Range("N40").Value = -10172
$ws.Range("H46").Value = 1350
$ws.Range("I46").Value = 1200
$ws.Range("J46").Value = 1500
$ws.Range("K46").Value = 1200
$ws.Range("L46").Value = 1500
$ws.Range("M46").Value = -1012
$ws.Range("N46").Value = -1876
$ws.Range("H55").Value = 10778.333
$ws.Range("I55").Value = 847.75
$ws.Range("J55").Value = 22127.572
$ws.Range("K55").Value = 847.75
$ws.Range("L55").Value = 22127.572
$ws.Range("M55").Value = -674.75
$ws.Range("N55").Value = -22473.572
$ws.Range("H61").Value = 1911.1666
$ws.Range("I61").Value = 1826.2222
$ws.Range("J61").Value = 2166
$ws.Range("K61").Value = 1826.2222
$ws.Range("L61").Value = 2166
$ws.Range("M61").Value = -1624.2222
$ws.Range("N61").Value = -2570
$ws.Range("H82").Value = 2599.8
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 2599.8
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H93").Value = 58824976
$ws.Range("I93").Value = 934.7857
$ws.Range("K93").Value = 934.7857
$ws.Range("M93").Value = 313.2143
$ws.Range("H98").Value = 54249.25
$ws.Range("J98").Value = 54249.25
$ws.Range("L98").Value = 54249.25
$ws.Range("N98").Value = -60239.25
$ws.Range("H104").Value = 47473
$ws.Range("J104").Value = 47473
$ws.Range("L104").Value = 47473
$ws.Range("N104").Value = -54461
$ws.Range("H106").Value = 25141
$ws.Range("J106").Value = 25141
$ws.Range("L106").Value = 25141
$ws.Range("N106").Value = -27665
$ws.Range("H113").Value = 1911.1666
$ws.Range("I113").Value = 1826.2222
$ws.Range("J113").Value = 2166
$ws.Range("K113").Value = 1826.2222
$ws.Range("L113").Value = 2166
$ws.Range("M113").Value = 343.7778000000001
$ws.Range("N113").Value = -6506
$ws.Range("H122").Value = 9030.444
$ws.Range("I122").Value = 8630.166999999999
$ws.Range("K122").Value = 25890.501
$ws.Range("M122").Value = -23440.501
$ws.Range("H126").Value = 8749.154
$ws.Range("I126").Value = 9118.849
$ws.Range("J126").Value = 6715.8335
$ws.Range("K126").Value = 27356.547
$ws.Range("L126").Value = 20147.5005
$ws.Range("M126").Value = -24886.547
$ws.Range("N126").Value = -25087.5005
$ws.Range("H132").Value = 6874.8335
$ws.Range("I132").Value = 6368.6553
$ws.Range("J132").Value = 7462
$ws.Range("K132").Value = 19105.9659
$ws.Range("L132").Value = 22386
$ws.Range("M132").Value = -16575.9659
$ws.Range("N132").Value = -27446

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 14799.5
$ws.Range("J24").Value = 14799.5
$ws.Range("L24").Value = 14799.5
$ws.Range("N24").Value = -15259.5
$ws.Range("H81").Value = 3976.2856
$ws.Range("I81").Value = 3122.4
$ws.Range("J81").Value = 6111
$ws.Range("K81").Value = 6244.8
$ws.Range("L81").Value = 12222
$ws.Range("M81").Value = -5183.8
$ws.Range("N81").Value = -14344
$ws.Range("H84").Value = 3976.2856
$ws.Range("I84").Value = 3122.4
$ws.Range("J84").Value = 6111
$ws.Range("K84").Value = 31224
$ws.Range("L84").Value = 61110
$ws.Range("M84").Value = -25920
$ws.Range("N84").Value = -71718
$ws.Range("H104").Value = 41997.5
$ws.Range("J104").Value = 41997.5
$ws.Range("L104").Value = 41997.5
$ws.Range("N104").Value = -48985.5
$ws.Range("H107").Value = 842.1429000000001
$ws.Range("I107").Value = 891.25
$ws.Range("K107").Value = 2673.75
$ws.Range("M107").Value = -753.75
$ws.Range("H122").Value = 1727.2858
$ws.Range("I122").Value = 1738.6
$ws.Range("K122").Value = 5215.799999999999
$ws.Range("M122").Value = -2765.799999999999
$ws.Range("H126").Value = 2381.8333
$ws.Range("I126").Value = 2258.4
$ws.Range("J126").Value = 2999
$ws.Range("K126").Value = 6775.200000000001
$ws.Range("L126").Value = 8997
$ws.Range("M126").Value = -4305.200000000001
$ws.Range("N126").Value = -13937
$ws.Range("H132").Value = 2505.7856
$ws.Range("I132").Value = 2523.7
$ws.Range("J132").Value = 2495.8333
$ws.Range("K132").Value = 7571.099999999999
$ws.Range("L132").Value = 7487.499899999999
$ws.Range("M132").Value = -5041.099999999999
$ws.Range("N132").Value = -12547.4999
$ws.Range("H136").Value = 318431.75
$ws.Range("I136").Value = 339536.7
$ws.Range("J136").Value = 1857.5
$ws.Range("K136").Value = 1018610.1
$ws.Range("L136").Value = 5572.5
$ws.Range("M136").Value = -1016060.1
$ws.Range("N136").Value = -10672.5
